$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Update the "Описание" (Description) column to reflect that the checks
# have now been performed ("Выполнена ..." prefix).
$ws.Range("B3").Value = "Выполнена проверка того, что видеоролик воспроизводится корректно"
$ws.Range("B4").Value = "Выполнена проверка на то, что кнопка работает согласно ожидаемому сценарию"
$ws.Range("B5").Value = "Выполнена проверка добавления комментария"
$ws.Range("B6").Value = "Выполнена проверка, что подписка производится корректно"
$ws.Range("B7").Value = "Выполнена проверка того, что осуществляется переход на страницу https://mail.ru/"

# Row 4 now needs more vertical space for the longer wrapped text.
$ws.Rows.Item(4).RowHeight = 60

# Move the active selection.
$ws.Range("J11").Select()
